$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 962.8
$ws.Range("I33").Value = 646.8889
$ws.Range("K33").Value = 646.8889
$ws.Range("M33").Value = -417.8889
# Row 69
$ws.Range("H69").Value = 21185
$ws.Range("I69").Value = 9000
$ws.Range("K69").Value = 27000
$ws.Range("M69").Value = -26126
# Row 70
$ws.Range("H70").Value = 877.1111
$ws.Range("I70").Value = 549.5
$ws.Range("J70").Value = 970.7143
$ws.Range("K70").Value = 1648.5
$ws.Range("L70").Value = 2912.1429
$ws.Range("M70").Value = -1378.5
$ws.Range("N70").Value = -3452.1429
# Row 72
$ws.Range("H72").Value = 21185
$ws.Range("I72").Value = 9000
$ws.Range("K72").Value = 81000
$ws.Range("M72").Value = -76632
# Row 73
$ws.Range("H73").Value = 877.1111
$ws.Range("I73").Value = 549.5
$ws.Range("J73").Value = 970.7143
$ws.Range("K73").Value = 1648.5
$ws.Range("L73").Value = 2912.1429
$ws.Range("M73").Value = -712.5
$ws.Range("N73").Value = -4784.1429
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = 0
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = 0
# Row 112
$ws.Range("H112").Value = 2332.2104
$ws.Range("I112").Value = 1345
$ws.Range("J112").Value = 2448.353
$ws.Range("K112").Value = 4035
$ws.Range("L112").Value = 7345.059
$ws.Range("M112").Value = -2927
$ws.Range("N112").Value = -9561.059000000001
# Row 127
$ws.Range("H127").Value = 699.6
# Row 132
$ws.Range("H132").Value = 5478.359
$ws.Range("I132").Value = 5469.697
$ws.Range("K132").Value = 16409.091
$ws.Range("M132").Value = -13879.091
# Row 137
$ws.Range("H137").Value = 12929.105
$ws.Range("I137").Value = 15443.533
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 46330.599
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -43780.599
$ws.Range("N137").Value = -15600
# Row 141
$ws.Range("H141").Value = 16943.75
$ws.Range("I141").Value = 40000
$ws.Range("J141").Value = 9258.333000000001
$ws.Range("K141").Value = 120000
$ws.Range("L141").Value = 27774.999
$ws.Range("M141").Value = -114820
$ws.Range("N141").Value = -38134.999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 450
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
# Row 5
$ws.Range("H5").Value = 702.5
$ws.Range("I5").Value = 716.3333
$ws.Range("K5").Value = 716.3333
$ws.Range("M5").Value = -604.3333
# Row 32
$ws.Range("H32").Value = 5496.4443
$ws.Range("I32").Value = 5231.5586
$ws.Range("K32").Value = 5231.5586
$ws.Range("M32").Value = -4944.5586
# Row 74
$ws.Range("H74").Value = 1587.8438
$ws.Range("I74").Value = 907.85187
$ws.Range("J74").Value = 5259.8
$ws.Range("K74").Value = 907.85187
$ws.Range("L74").Value = 5259.8
$ws.Range("M74").Value = -33.85186999999996
$ws.Range("N74").Value = -7007.8
# Row 77
$ws.Range("H77").Value = 1587.8438
$ws.Range("I77").Value = 907.85187
$ws.Range("J77").Value = 5259.8
$ws.Range("K77").Value = 4539.25935
$ws.Range("L77").Value = 26299
$ws.Range("M77").Value = -171.2593500000003
$ws.Range("N77").Value = -35035

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 702.5
$ws.Range("I4").Value = 716.3333
$ws.Range("K4").Value = 716.3333
$ws.Range("M4").Value = -601.3333
# Row 20
$ws.Range("H20").Value = 2835.818
$ws.Range("I20").Value = 1724
$ws.Range("J20").Value = 4781.5
$ws.Range("K20").Value = 1724
$ws.Range("L20").Value = 4781.5
$ws.Range("M20").Value = -1477
$ws.Range("N20").Value = -5275.5
# Row 134
$ws.Range("H134").Value = 7273.0435
$ws.Range("I134").Value = 8398
$ws.Range("K134").Value = 25194
$ws.Range("M134").Value = -22659
# Row 141
$ws.Range("H141").Value = 150000
$ws.Range("J141").Value = 150000
$ws.Range("L141").Value = 150000
$ws.Range("N141").Value = -160360

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7188.2964
$ws.Range("I31").Value = 7211.875
$ws.Range("J31").Value = 6999.6665
$ws.Range("K31").Value = 7211.875
$ws.Range("L31").Value = 6999.6665
$ws.Range("M31").Value = -6916.875
$ws.Range("N31").Value = -7589.6665
# Row 34
$ws.Range("H34").Value = 7188.2964
$ws.Range("I34").Value = 7211.875
$ws.Range("J34").Value = 6999.6665
$ws.Range("K34").Value = 7211.875
$ws.Range("L34").Value = 6999.6665
$ws.Range("M34").Value = -7009.875
$ws.Range("N34").Value = -7403.6665
# Row 132
$ws.Range("H132").Value = 2142.6365
$ws.Range("I132").Value = 2061.6667
$ws.Range("J132").Value = 2507
$ws.Range("K132").Value = 6185.000100000001
$ws.Range("L132").Value = 7521
$ws.Range("M132").Value = -3655.000100000001
$ws.Range("N132").Value = -12581
# Row 134
$ws.Range("H134").Value = 2040.3462
$ws.Range("I134").Value = 1958.7391
$ws.Range("K134").Value = 5876.2173
$ws.Range("M134").Value = -3341.2173

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 1701.091
$ws.Range("I38").Value = 356.5
$ws.Range("J38").Value = 2469.4285
$ws.Range("K38").Value = 1069.5
$ws.Range("L38").Value = 7408.2855
$ws.Range("M38").Value = -722.5
$ws.Range("N38").Value = -8102.2855
# Row 50
$ws.Range("H50").Value = 1677.3077
$ws.Range("I50").Value = 1225
$ws.Range("J50").Value = 1878.3334
$ws.Range("K50").Value = 3675
$ws.Range("L50").Value = 5635.0002
$ws.Range("M50").Value = -3194
$ws.Range("N50").Value = -6597.0002
# Row 53
$ws.Range("H53").Value = 1677.3077
$ws.Range("I53").Value = 1225
$ws.Range("J53").Value = 1878.3334
$ws.Range("K53").Value = 3675
$ws.Range("L53").Value = 5635.0002
$ws.Range("M53").Value = -3194
$ws.Range("N53").Value = -6597.0002
# Row 60
$ws.Range("H60").Value = 1070.8572
$ws.Range("I60").Value = 465.83334
$ws.Range("J60").Value = 1524.625
$ws.Range("K60").Value = 1397.50002
$ws.Range("L60").Value = 4573.875
$ws.Range("M60").Value = -1146.50002
$ws.Range("N60").Value = -5075.875
# Row 109
$ws.Range("H109").Value = 2600.5
$ws.Range("J109").Value = 3499.5
$ws.Range("L109").Value = 10498.5
$ws.Range("N109").Value = -12578.5
# Row 131
$ws.Range("H131").Value = 1618.8636
$ws.Range("J131").Value = 1634.9025
$ws.Range("L131").Value = 4904.7075
$ws.Range("N131").Value = -14984.7075

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 11753.875
$ws.Range("I70").Value = 11376
$ws.Range("J70").Value = 12131.75
$ws.Range("K70").Value = 11376
$ws.Range("L70").Value = 12131.75
$ws.Range("M70").Value = -11106
$ws.Range("N70").Value = -12671.75
# Row 73
$ws.Range("H73").Value = 11753.875
$ws.Range("I73").Value = 11376
$ws.Range("J73").Value = 12131.75
$ws.Range("K73").Value = 11376
$ws.Range("L73").Value = 12131.75
$ws.Range("M73").Value = -10440
$ws.Range("N73").Value = -14003.75
# Row 113
$ws.Range("H113").Value = 11433.667
$ws.Range("I113").Value = 23288.8
$ws.Range("J113").Value = 2965.7144
$ws.Range("K113").Value = 23288.8
$ws.Range("L113").Value = 2965.7144
$ws.Range("M113").Value = -21118.8
$ws.Range("N113").Value = -7305.7144
# Row 126
$ws.Range("H126").Value = 8709.0625
$ws.Range("I126").Value = 10756.429
$ws.Range("J126").Value = 7116.6665
$ws.Range("K126").Value = 32269.287
$ws.Range("L126").Value = 21349.9995
$ws.Range("M126").Value = -29799.287
$ws.Range("N126").Value = -26289.9995
# Row 132
$ws.Range("H132").Value = 2492.8
$ws.Range("I132").Value = 2277.1428
$ws.Range("K132").Value = 6831.428400000001
$ws.Range("M132").Value = -4301.428400000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 29209.572
$ws.Range("J7").Value = 7999.4287
$ws.Range("L7").Value = 7999.4287
$ws.Range("N7").Value = -8223.4287
# Row 40
$ws.Range("H40").Value = 28799.428
$ws.Range("I40").Value = 41009.816
$ws.Range("J40").Value = 15368
$ws.Range("K40").Value = 41009.816
$ws.Range("L40").Value = 15368
$ws.Range("M40").Value = -40873.816
$ws.Range("N40").Value = -15640
# Row 42
$ws.Range("H42").Value = 87400
$ws.Range("I42").Value = 87400
$ws.Range("K42").Value = 87400
$ws.Range("M42").Value = -86837
# Row 49
$ws.Range("H49").Value = 87400
$ws.Range("I49").Value = 87400
$ws.Range("K49").Value = 87400
$ws.Range("M49").Value = -87253
# Row 53
$ws.Range("H53").Value = 22000
$ws.Range("J53").Value = 22000
$ws.Range("L53").Value = 22000
$ws.Range("N53").Value = -23036
# Row 100
$ws.Range("H100").Value = 4011.8667
$ws.Range("I100").Value = 2023.875
$ws.Range("J100").Value = 6283.857
$ws.Range("K100").Value = 2023.875
$ws.Range("L100").Value = 6283.857
$ws.Range("M100").Value = -1482.875
$ws.Range("N100").Value = -7365.857
# Row 126
$ws.Range("H126").Value = 29209.572
$ws.Range("J126").Value = 7999.4287
$ws.Range("L126").Value = 23998.2861
$ws.Range("N126").Value = -28938.2861
# Row 132
$ws.Range("H132").Value = 681001.5600000001
$ws.Range("I132").Value = 934708.75
$ws.Range("J132").Value = 4449
$ws.Range("K132").Value = 2804126.25
$ws.Range("L132").Value = 13347
$ws.Range("M132").Value = -2801596.25
$ws.Range("N132").Value = -18407
# Row 136
$ws.Range("H136").Value = 5684.5386
$ws.Range("I136").Value = 3729.1428
$ws.Range("K136").Value = 11187.4284
$ws.Range("M136").Value = -8637.428400000001
